$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Roraima"
$ws.Range("B2").Value = "Variação 2023/2022"
$ws.Range("C2").Value = 22.01060354937823

# Row 3
$ws.Range("A3").Value = "Tocantins"
$ws.Range("B3").Value = "Variação 2023/2022"
$ws.Range("C3").Value = 17.1014622367303

# Row 4
$ws.Range("A4").Value = "Piauí"
$ws.Range("B4").Value = "Variação 2023/2022"
$ws.Range("C4").Value = 14.76341704251287

# Row 5
$ws.Range("A5").Value = "Amapá"
$ws.Range("B5").Value = "Variação 2023/2022"
$ws.Range("C5").Value = 13.26378113223738

# Row 6
$ws.Range("A6").Value = "Mato Grosso"
$ws.Range("B6").Value = "Variação 2023/2022"
$ws.Range("C6").Value = 11.47959452658822

# Row 7
$ws.Range("A7").Value = "Acre"
$ws.Range("B7").Value = "Variação 2023/2022"
$ws.Range("C7").Value = 11.41224342862286

# Row 8
$ws.Range("B8").Value = "Variação 2023/2022"
$ws.Range("C8").Value = 8.383137353614533
$ws.Range("D8").Value = "7º"

# Row 9
$ws.Range("B9").Value = "Variação 2023/2022"
$ws.Range("C9").Value = 5.617270641488692

# Row 10
$ws.Range("B10").Value = "Variação 2023/2022"
$ws.Range("C10").Value = 3.628508523747587
